$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the INVI (column E) values for the affected row blocks.
$ws.Range("E2:E21").Value = "John Doe0"
$ws.Range("E42:E61").Value = "John Doe10"
$ws.Range("E62:E81").Value = "John Doe15"
$ws.Range("E82:E101").Value = "John Doe20"
$ws.Range("E102:E121").Value = "John Doe25"
$ws.Range("E122:E141").Value = "John Doe30"
$ws.Range("E142:E161").Value = "John Doe35"
$ws.Range("E162:E181").Value = "John Doe40"
$ws.Range("E182:E201").Value = "John Doe45"

# Rows 202-393 no longer have an INVI value assigned.
$ws.Range("E202:E393").ClearContents()
